$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. GrthDCmp block: collapse 7 old rows (117-123) into 6 new rows (117-122) ---
# Delete one row from the old block (row 120) so rows below shift up by one.
# (This single-row deletion also cascades the trailing "ISO50001" row from
#  281 down to 280, which is the entirety of the apparent "row 281 removed"
#  portion of the diff -- the content itself is unchanged, only its row
#  number shifts.)
$ws.Rows("120:120").Delete()

# Now rows 117-122 hold (post-shift): 117,118,119 = old 117,118,119 ; 120,121,122 = old 121,122,123
# Overwrite them with the new content.
$ws.Range("A117").Value = "GrthDCmpNC"
$ws.Range("B117").Value = "Final consumption expenditure (% of GDP growth)"
$ws.Range("C117").Value = "Dépense de consommation finale (% de croissance du PIB)"

$ws.Range("A118").Value = "GrthDCmpNFI"
$ws.Range("B118").Value = "Gross fixed capital formation (% of GDP growth)"
$ws.Range("C118").Value = "Formation brute de capital fixe (% de croissance du PIB)"

$ws.Range("A119").Value = "GrthDCmpNINV"
$ws.Range("B119").Value = "Changes in inventories (% of GDP growth)"
$ws.Range("C119").Value = "Variations des stocks (% de croissance du PIB)"

$ws.Range("A120").Value = "GrthDCmpNX"
$ws.Range("B120").Value = "Exports of goods and services (% of GDP growth)"
$ws.Range("C120").Value = "Exportations de biens et services (% de croissance du PIB)"

$ws.Range("A121").Value = "GrthDCmpNM"
$ws.Range("B121").Value = "Imports of goods and services (% of GDP growth)"
$ws.Range("C121").Value = "Importations de biens et services (% de croissance du PIB)"

$ws.Range("A122").Value = "GrthDCmpNFB"
$ws.Range("B122").Value = "Foreign balance (% of GDP growth)"
$ws.Range("C122").Value = "Balance extérieure (% de croissance du PIB)"

# --- 3. Update the _xlnm._FilterDatabase defined name range to the new extent ---
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=AfDD_DDAf_2022_TabIndsCodes!`$A`$1:`$C`$280"
    }
}

# --- 4. Selection state shown in the diff ---
$ws.Range("B122").Select()
